# Applies the commit: inserts 3 new price rows (Melón / Tuna, 05-12-2022)
# into the "Femacal de La Calera" subset sheet, above the old row 374,
# shifting all subsequent rows down by 3 (old 374:474 -> new 377:477).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at position 374 - this pushes the existing
# rows 374:474 down to 377:477, preserving their values/formatting.
$ws.Rows("374:376").Insert()

# --- New row 374: Tuna / Extra ---
$ws.Range("A374").Value2 = 3
$ws.Range("B374").Value2 = "Femacal de La Calera"
$ws.Range("C374").Value2 = "Coquimbo"
$ws.Range("D374").Value2 = 44900
$ws.Range("E374").Value2 = 5
$ws.Range("F374").Value2 = 100112027
$ws.Range("G374").Value2 = "Melón"
$ws.Range("H374").Value2 = "Tuna"
$ws.Range("I374").Value2 = "Extra"
$ws.Range("J374").Value2 = 700
$ws.Range("K374").Value2 = 2000
$ws.Range("L374").Value2 = 2000
$ws.Range("M374").Value2 = 2000
$ws.Range("N374").Value2 = "$/unidad"
$ws.Range("O374").Value2 = "Provincia de Quillota"
$ws.Range("P374").Value2 = 2000
$ws.Range("Q374").Value2 = 1
$ws.Range("R374").Value2 = "Hortaliza"

# --- New row 375: Tuna / Primera ---
$ws.Range("A375").Value2 = 3
$ws.Range("B375").Value2 = "Femacal de La Calera"
$ws.Range("C375").Value2 = "Coquimbo"
$ws.Range("D375").Value2 = 44900
$ws.Range("E375").Value2 = 5
$ws.Range("F375").Value2 = 100112027
$ws.Range("G375").Value2 = "Melón"
$ws.Range("H375").Value2 = "Tuna"
$ws.Range("I375").Value2 = "Primera"
$ws.Range("J375").Value2 = 700
$ws.Range("K375").Value2 = 1500
$ws.Range("L375").Value2 = 1500
$ws.Range("M375").Value2 = 1500
$ws.Range("N375").Value2 = "$/unidad"
$ws.Range("O375").Value2 = "Provincia de Quillota"
$ws.Range("P375").Value2 = 1500
$ws.Range("Q375").Value2 = 1
$ws.Range("R375").Value2 = "Hortaliza"

# --- New row 376: Tuna / Segunda ---
$ws.Range("A376").Value2 = 3
$ws.Range("B376").Value2 = "Femacal de La Calera"
$ws.Range("C376").Value2 = "Coquimbo"
$ws.Range("D376").Value2 = 44900
$ws.Range("E376").Value2 = 5
$ws.Range("F376").Value2 = 100112027
$ws.Range("G376").Value2 = "Melón"
$ws.Range("H376").Value2 = "Tuna"
$ws.Range("I376").Value2 = "Segunda"
$ws.Range("J376").Value2 = 650
$ws.Range("K376").Value2 = 1000
$ws.Range("L376").Value2 = 1000
$ws.Range("M376").Value2 = 1000
$ws.Range("N376").Value2 = "$/unidad"
$ws.Range("O376").Value2 = "Provincia de Quillota"
$ws.Range("P376").Value2 = 1000
$ws.Range("Q376").Value2 = 1
$ws.Range("R376").Value2 = "Hortaliza"

Write-Host "Done: inserted rows 374-376, dimension should now be A1:R477"
